$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 298 (the weekly update adds a new price-report date,
# pushing every existing record for this market/product down by 2 rows).
$ws.Rows("298:299").Insert()

# New row 298: Lechuga - Conconina(o), 2021-09-09, Región Metropolitana
$ws.Cells.Item(298, 1).Value = 4
$ws.Cells.Item(298, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(298, 3).Value = "Los Lagos"
$ws.Cells.Item(298, 4).Value = 44448
$ws.Cells.Item(298, 5).Value = 10
$ws.Cells.Item(298, 6).Value = 100112033
$ws.Cells.Item(298, 7).Value = "Lechuga"
$ws.Cells.Item(298, 8).Value = "Conconina(o)"
$ws.Cells.Item(298, 9).Value = "Primera"
$ws.Cells.Item(298, 10).Value = 110
$ws.Cells.Item(298, 11).Value = 11000
$ws.Cells.Item(298, 12).Value = 11000
$ws.Cells.Item(298, 13).Value = 11000
$ws.Cells.Item(298, 14).Value = "`$/caja 10 unidades"
$ws.Cells.Item(298, 15).Value = "Región Metropolitana"
$ws.Cells.Item(298, 16).Value = 1100
$ws.Cells.Item(298, 17).Value = 10
$ws.Cells.Item(298, 18).Value = "Hortaliza"

# New row 299: Lechuga - Escarola, 2021-09-09, Región de Coquimbo
$ws.Cells.Item(299, 1).Value = 4
$ws.Cells.Item(299, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(299, 3).Value = "Los Lagos"
$ws.Cells.Item(299, 4).Value = 44448
$ws.Cells.Item(299, 5).Value = 10
$ws.Cells.Item(299, 6).Value = 100112033
$ws.Cells.Item(299, 7).Value = "Lechuga"
$ws.Cells.Item(299, 8).Value = "Escarola"
$ws.Cells.Item(299, 9).Value = "Primera"
$ws.Cells.Item(299, 10).Value = 250
$ws.Cells.Item(299, 11).Value = 13000
$ws.Cells.Item(299, 12).Value = 13500
$ws.Cells.Item(299, 13).Value = 13300
$ws.Cells.Item(299, 14).Value = "`$/caja 15 unidades"
$ws.Cells.Item(299, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(299, 16).Value = 887
$ws.Cells.Item(299, 17).Value = 15
$ws.Cells.Item(299, 18).Value = "Hortaliza"
